$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new payment entry: label in H17, amount in H18
$ws.Range("H17").Value = "Check 7/21/16"
$ws.Range("H18").Value = 242.38

# Update the selected cell to reflect the new active selection
$ws.Range("H19").Select()
